$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two data rows (gs, dice) that fell off the bottom of the
# new, shorter table. This also shrinks the sheet dimension from F13 to F11.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(12).Delete()

# CERScore recomputed as the inverse of distance -> new values for rows 3-11.
# Row 2 (KAOGExp) is unchanged.

$ws.Range("A3").Value = "face-knn"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 9.835418325367966
$ws.Range("D3").Value = 21.95489729689547
$ws.Range("E3").Value = 26.503318857581
$ws.Range("F3").Value = 75.60882951000001

$ws.Range("A4").Value = "dice"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 40.91807604500001
$ws.Range("D4").Value = 76.92589078860948
$ws.Range("E4").Value = 215.1916095887813
$ws.Range("F4").Value = 111.2721740379802

$ws.Range("A5").Value = "clue"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 10.82941470457071
$ws.Range("D5").Value = 34.294330352762
$ws.Range("E5").Value = 61.15496029938544
$ws.Range("F5").Value = 96.242649608826

$ws.Range("A6").Value = "cruds"
$ws.Range("B6").Value = 0.45
$ws.Range("C6").Value = 2.502678916033189
$ws.Range("D6").Value = 5.233592245180635
$ws.Range("E6").Value = 6.850932417469378
$ws.Range("F6").Value = 21.91949592025048

$ws.Range("A7").Value = "cem"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 10.32785084621429
$ws.Range("D7").Value = 490.2476918543464
$ws.Range("E7").Value = 25698.92363656488
$ws.Range("F7").Value = 494.5480324535401

$ws.Range("A8").Value = "wachter"
$ws.Range("B8").Value = 0.5600000000000001
$ws.Range("C8").Value = 2.801242255
$ws.Range("D8").Value = 77.04564214242079
$ws.Range("E8").Value = 5008.075595881454
$ws.Range("F8").Value = 414.1865435513668

$ws.Range("A9").Value = "face-epsilon"
$ws.Range("B9").Value = 0.95
$ws.Range("C9").Value = 8.808593235964647
$ws.Range("D9").Value = 18.9142653180085
$ws.Range("E9").Value = 22.28621378253619
$ws.Range("F9").Value = 64.9514995

$ws.Range("A10").Value = "gs"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 10.56844646036905
$ws.Range("D10").Value = 128.4959148244709
$ws.Range("E10").Value = 5405.51302451188
$ws.Range("F10").Value = 305.8624585226183

$ws.Range("A11").Value = "ar"
$ws.Range("B11").Value = 0.26
$ws.Range("C11").Value = 2.91803754
$ws.Range("D11").Value = 3.965654871881759
$ws.Range("E11").Value = 3.979414772917486
$ws.Range("F11").Value = 4.96790886
